# Update 2015年-2020年 (rows 4-9) with higher-precision values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 22.7129126794
$ws.Range("C4").Value = 119.931238329
$ws.Range("D4").Value = 36.9044133806
$ws.Range("E4").Value = 45.7188785068
$ws.Range("F4").Value = 42.2155958726
$ws.Range("G4").Value = 86.4123545881
$ws.Range("H4").Value = 71.1645900932
$ws.Range("I4").Value = 20.3921549472
$ws.Range("J4").Value = 88.9919479417
$ws.Range("K4").Value = 47.6269806615
$ws.Range("L4").Value = 224.828289542
$ws.Range("M4").Value = 81.5241847911
$ws.Range("N4").Value = 55.4954585058

$ws.Range("B5").Value = 27.6897
$ws.Range("C5").Value = 120.7848
$ws.Range("D5").Value = 38.4258
$ws.Range("E5").Value = 48.6116
$ws.Range("F5").Value = 39.9532
$ws.Range("G5").Value = 89.8212
$ws.Range("H5").Value = 76.1649
$ws.Range("I5").Value = 17.7144
$ws.Range("J5").Value = 93.4693
$ws.Range("K5").Value = 53.1609
$ws.Range("L5").Value = 235.4053
$ws.Range("M5").Value = 90.8779
$ws.Range("N5").Value = 57.5458

$ws.Range("B6").Value = 29.7344154567
$ws.Range("C6").Value = 122.1810080929
$ws.Range("D6").Value = 39.9859459599
$ws.Range("E6").Value = 50.9600943565
$ws.Range("F6").Value = 39.303497272
$ws.Range("G6").Value = 91.6596050454
$ws.Range("H6").Value = 78.6309920766
$ws.Range("I6").Value = 18.3540340257
$ws.Range("J6").Value = 95.3359341247
$ws.Range("K6").Value = 56.4848135541
$ws.Range("L6").Value = 239.9520319968
$ws.Range("M6").Value = 96.1469674323
$ws.Range("N6").Value = 58.7364643101

$ws.Range("B7").Value = 32.9989570145
$ws.Range("C7").Value = 119.2620095307
$ws.Range("D7").Value = 39.2066561147
$ws.Range("E7").Value = 56.3998196869
$ws.Range("F7").Value = 35.7102033761
$ws.Range("G7").Value = 93.7802110135
$ws.Range("H7").Value = 84.997862906
$ws.Range("I7").Value = 12.6331033588
$ws.Range("J7").Value = 98.7604397345
$ws.Range("K7").Value = 59.2116994735
$ws.Range("L7").Value = 249.0600017047
$ws.Range("M7").Value = 109.2980589505
$ws.Range("N7").Value = 53.3637600073

$ws.Range("B8").Value = 35.3122365057
$ws.Range("C8").Value = 120.6009669851
$ws.Range("D8").Value = 40.1274177579
$ws.Range("E8").Value = 59.3165244781
$ws.Range("F8").Value = 34.1528952698
$ws.Range("G8").Value = 95.9612615623
$ws.Range("H8").Value = 86.9439596707
$ws.Range("I8").Value = 12.2222278616
$ws.Range("J8").Value = 100.8668101275
$ws.Range("K8").Value = 63.9379587444
$ws.Range("L8").Value = 253.2403151615
$ws.Range("M8").Value = 115.643282159
$ws.Range("N8").Value = 53.2371923715

$ws.Range("B9").Value = 37.1401929058467
$ws.Range("C9").Value = 120.797903702809
$ws.Range("D9").Value = 41.0047492149819
$ws.Range("E9").Value = 60.8557657646377
$ws.Range("F9").Value = 33.1081645026222
$ws.Range("G9").Value = 96.74755570737879
$ws.Range("H9").Value = 90.4325240756464
$ws.Range("I9").Value = 12.1273951810816
$ws.Range("J9").Value = 101.834619138211
$ws.Range("K9").Value = 66.6891953518998
$ws.Range("L9").Value = 253.840843772832
$ws.Range("M9").Value = 117.738014531804
$ws.Range("N9").Value = 54.1889584463376

# Add new row 10: 2021年
$ws.Range("A10").Value = "2021年"
$ws.Range("B10").Value = 41.8
$ws.Range("C10").Value = 118.7
$ws.Range("D10").Value = 41.5
$ws.Range("E10").Value = 63.1
$ws.Range("F10").Value = 31.5
$ws.Range("G10").Value = 98.7
$ws.Range("H10").Value = 89.59999999999999
$ws.Range("I10").Value = 8.1
$ws.Range("J10").Value = 103.9
$ws.Range("K10").Value = 73.8
$ws.Range("L10").Value = 259.1
$ws.Range("M10").Value = 131.2
$ws.Range("N10").Value = 47

# Copy label-column formatting (bold, centered, bordered) from row 9 to row 10
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats

# Add new row 11: 2022年 (some indicators not yet reported -> left blank)
$ws.Range("A11").Value = "2022年"
$ws.Range("B11").Value = 43.5419742834177
$ws.Range("C11").Value = 118.885001318792
$ws.Range("F11").Value = 30.8052879397484
$ws.Range("G11").Value = 99.04344805616149
$ws.Range("J11").Value = 104.170812652236
$ws.Range("L11").Value = 259.351694347994
$ws.Range("M11").Value = 133.898697356459
$ws.Range("N11").Value = 47.4604270421034

# Copy label-column formatting (bold, centered, bordered) from row 9 to row 11
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

Write-Host "Edit complete"
